# Update stock metrics (current price, change rate, ROE-ish col I, Stochastic %K/%D) per latest data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24100
$ws.Range("D2").Value = 0.019
$ws.Range("I2").Value = 4.15
$ws.Range("J2").Value = 98
$ws.Range("K2").Value = 98
$ws.Range("C3").Value = 100800
$ws.Range("D3").Value = 0.007
$ws.Range("I3").Value = 6.45
$ws.Range("C4").Value = 473000
$ws.Range("D4").Value = 0.0042
$ws.Range("I4").Value = 4.02
$ws.Range("J4").Value = 97
$ws.Range("K4").Value = 97
$ws.Range("C5").Value = 31900
$ws.Range("D5").Value = 0.0095
$ws.Range("I5").Value = 6.27
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("C6").Value = 31000
$ws.Range("D6").Value = -0.0112
$ws.Range("I6").Value = 3.87
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 80
$ws.Range("C7").Value = 31250
$ws.Range("D7").Value = 0.0434
$ws.Range("I7").Value = 3.84
$ws.Range("J7").Value = 99
$ws.Range("K7").Value = 99
$ws.Range("C8").Value = 11150
$ws.Range("D8").Value = -0.0168
$ws.Range("I8").Value = 4.62
$ws.Range("J8").Value = 96
$ws.Range("K8").Value = 96
$ws.Range("C9").Value = 76000
$ws.Range("D9").Value = -0.0065
$ws.Range("I9").Value = 3.95
$ws.Range("J9").Value = 62
$ws.Range("K9").Value = 62
$ws.Range("C10").Value = 208000
$ws.Range("C11").Value = 138800
$ws.Range("D11").Value = 0.0154
$ws.Range("I11").Value = 4.9
$ws.Range("J11").Value = 96
$ws.Range("K11").Value = 96
$ws.Range("C12").Value = 22650
$ws.Range("D12").Value = 0.0157
$ws.Range("I12").Value = 4.19
$ws.Range("J12").Value = 98
$ws.Range("K12").Value = 98
$ws.Range("C13").Value = 77300
$ws.Range("C14").Value = 55700
$ws.Range("D14").Value = 0.0072
$ws.Range("I14").Value = 6.36
$ws.Range("J14").Value = 72
$ws.Range("K14").Value = 72
$ws.Range("C15").Value = 87400
$ws.Range("D15").Value = 0.0058
$ws.Range("I15").Value = 6.29
$ws.Range("J15").Value = 94
$ws.Range("K15").Value = 94
$ws.Range("C16").Value = 21000
$ws.Range("D16").Value = 0.037
$ws.Range("I16").Value = 5.07
$ws.Range("C17").Value = 56600
$ws.Range("D17").Value = 0.0107
$ws.Range("I17").Value = 4.95
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 97
$ws.Range("C18").Value = 21650
$ws.Range("D18").Value = -0.0023
$ws.Range("I18").Value = 5.68
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = 51
$ws.Range("C19").Value = 58100
$ws.Range("D19").Value = 0.0338
$ws.Range("I19").Value = 3.44
$ws.Range("C20").Value = 14990
$ws.Range("D20").Value = 0.0108
$ws.Range("I20").Value = 4.34
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 82
$ws.Range("C21").Value = 136400
$ws.Range("D21").Value = 0.0059
$ws.Range("I21").Value = 3.96
$ws.Range("J21").Value = 93
$ws.Range("K21").Value = 93
$ws.Range("C22").Value = 45500
$ws.Range("D22").Value = -0.0044
$ws.Range("I22").Value = 3.2
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = 54
$ws.Range("C23").Value = 69300
$ws.Range("D23").Value = -0.0086
$ws.Range("I23").Value = 3.12
$ws.Range("J23").Value = 94
$ws.Range("K23").Value = 94
$ws.Range("C24").Value = 52100
$ws.Range("D24").Value = -0.0133
$ws.Range("I24").Value = 5.18
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 80
$ws.Range("C25").Value = 93200
$ws.Range("D25").Value = 0.0108
$ws.Range("I25").Value = 3.86
$ws.Range("J25").Value = 96
$ws.Range("K25").Value = 96
$ws.Range("C26").Value = 117300
$ws.Range("D26").Value = 0.0043
$ws.Range("I26").Value = 2.71
$ws.Range("J26").Value = 95
$ws.Range("K26").Value = 95
$ws.Range("C27").Value = 15310
$ws.Range("D27").Value = 0.0153
$ws.Range("I27").Value = 4.25
$ws.Range("J27").Value = 98
$ws.Range("K27").Value = 98
$ws.Range("C28").Value = 14980
$ws.Range("D28").Value = -0.0013
$ws.Range("I28").Value = 3.34
$ws.Range("J28").Value = 99
$ws.Range("K28").Value = 99
$ws.Range("C29").Value = 25350
$ws.Range("D29").Value = 0.004
$ws.Range("I29").Value = 3.93
$ws.Range("C30").Value = 25850
$ws.Range("D30").Value = 0.0258
$ws.Range("I30").Value = 4.64
$ws.Range("J30").Value = 99
$ws.Range("K30").Value = 99

# D10 and D13 swap number format (0% <-> 0.00%) along with new values
$ws.Range("D10").NumberFormat = "0%"
$ws.Range("D10").Value = 0
$ws.Range("D13").NumberFormat = "0.00%"
$ws.Range("D13").Value = -0.0064

# Restore active cell selection
$ws.Range("E16").Select()
